# Generate Report for Handoff
#
# The localization run for 4f014fc9-9bb1-48cf-b7b7-ec71c489cbe8.md moved from
# "In Translation" to "Ready for handoff": refresh the status text and the
# associated handoff timestamps on the Overview roll-up sheet and on each
# per-locale (zh-cn / de-de) detail sheet, and widen the status/date columns
# slightly so the longer "Ready for handoff" label + refreshed timestamps
# aren't clipped.

$wb = $excel.ActiveWorkbook

# Column E/F (Overview) and column C (zh-cn, de-de) need to grow from ~13.41
# characters to ~17.22 characters. Excel's ColumnWidth setter only resolves
# to discrete ~1/6-character steps, so 16.3333... is the input that lands on
# the stored grid value closest to the target width.
$newStatusColWidth = 16.333333333333336

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"          # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"          # de-de status
$wsOverview.Range("G2").Value = "2016-08-18 10:43:19"        # Latest HO Xliff Generate Date
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"              # Status
$wsZhCn.Range("H2").Value = "2016-08-18 10:43:09"            # Latest Handoff Datetime
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"              # Status
$wsDeDe.Range("H2").Value = "2016-08-18 10:43:19"            # Latest Handoff Datetime
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
